$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "67.643.58"
Set-TextValue $ws.Range("E2") "  -3.11%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.522.50"
Set-TextValue $ws.Range("E3") "  -3.09%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.30%  "

# Row 5
Set-TextValue $ws.Range("D5") "608.89"
Set-TextValue $ws.Range("E5") "  -3.84%  "

# Row 6
Set-TextValue $ws.Range("D6") "150.33"
Set-TextValue $ws.Range("E6") "  -5.99%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.521.27"
Set-TextValue $ws.Range("E7") "  -3.02%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.14%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.487"
Set-TextValue $ws.Range("E9") "  -1.81%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.144"
Set-TextValue $ws.Range("E10") "  -3.67%  "

# Row 11
Set-TextValue $ws.Range("D11") "7.54"
Set-TextValue $ws.Range("E11") "  +3.74%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.430"
Set-TextValue $ws.Range("E12") "  -2.59%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000216"
Set-TextValue $ws.Range("E13") "  -6.41%  "

# Row 14
Set-TextValue $ws.Range("D14") "32.19"
Set-TextValue $ws.Range("E14") "  -3.83%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.083.98"
Set-TextValue $ws.Range("E15") "  -3.93%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.510.99"
Set-TextValue $ws.Range("E16") "  -3.48%  "

# Row 17
Set-TextValue $ws.Range("D17") "67.974.53"
Set-TextValue $ws.Range("E17") "  -2.27%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.117"
Set-TextValue $ws.Range("E18") "  -0.48%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.52"
Set-TextValue $ws.Range("E19") "  -2.44%  "

# Row 20
Set-TextValue $ws.Range("D20") "15.57"
Set-TextValue $ws.Range("E20") "  -3.14%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.00"
Set-TextValue $ws.Range("E21") "  -1.72%  "

# Row 22
Set-TextValue $ws.Range("D22") "453.37"
Set-TextValue $ws.Range("E22") "  -2.80%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.626"
Set-TextValue $ws.Range("E23") "  -3.02%  "

# Row 24
Set-TextValue $ws.Range("D24") "78.92"
Set-TextValue $ws.Range("E24") "  +0.08%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.648.17"
Set-TextValue $ws.Range("E25") "  -3.50%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.04%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0000123"
Set-TextValue $ws.Range("E27") "  -10.26%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "9.96"
Set-TextValue $ws.Range("E28") "  -7.31%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D29") "8.62"
Set-TextValue $ws.Range("E29") "  -7.59%  "

# Row 30
Set-TextValue $ws.Range("D30") "2.52"
Set-TextValue $ws.Range("E30") "  -5.20%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.66"
Set-TextValue $ws.Range("E31") "  -4.12%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -5.58%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.999"
Set-TextValue $ws.Range("E33") "  +0.06%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "6.25"
Set-TextValue $ws.Range("E34") "  -6.05%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D35") "25.76"
Set-TextValue $ws.Range("E35") "  -3.37%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.84"
Set-TextValue $ws.Range("E36") "  -6.87%  "

# Row 37
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D37") "3.489.90"
Set-TextValue $ws.Range("E37") "  -3.92%  "

# Row 38
Set-TextValue $ws.Range("D38") "8.02"
Set-TextValue $ws.Range("E38") "  -5.71%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +0.03%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.28"
Set-TextValue $ws.Range("E40") "  -6.78%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.995"
Set-TextValue $ws.Range("E41") "  -0.64%  "

# Row 42
Set-TextValue $ws.Range("D42") "175.58"
Set-TextValue $ws.Range("E42") "  -1.13%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.0900"
Set-TextValue $ws.Range("E43") "  -3.55%  "

# Row 44
Set-TextValue $ws.Range("D44") "5.44"
Set-TextValue $ws.Range("E44") "  -3.81%  "

# Row 45
Set-TextValue $ws.Range("D45") "31.05"
Set-TextValue $ws.Range("E45") "  -2.44%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.901"
Set-TextValue $ws.Range("E46") "  -1.80%  "

# Row 47
Set-TextValue $ws.Range("D47") "46.89"
Set-TextValue $ws.Range("E47") "  +0.77%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -6.54%  "

# Row 49
Set-TextValue $ws.Range("D49") "7.63"
Set-TextValue $ws.Range("E49") "  -2.57%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.49"
Set-TextValue $ws.Range("E50") "  -11.50%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.00"
Set-TextValue $ws.Range("E51") "  -3.44%  "
